# Trading update: 2026-02-17 12:40:08
# Appends a new open MarketMaking trade (row 40) to both the "All Trades"
# and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

function Add-TradeRow {
    param($ws, $row)

    $ws.Range("A$row").Value = 39

    # "2026-02-17" looks like a date, so a plain assignment would get
    # auto-converted into a date serial by Excel. Force literal text entry
    # (the way typing `'2026-02-17` into a cell would), then drop back to
    # the default "Normal" style so no stray number-format/quote-prefix
    # formatting is left behind on the cell.
    $ws.Range("B$row").Value = "'2026-02-17"
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").Value = "12:39:21"
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "UP"
    $ws.Range("F$row").Value = 0.97

    # Exit Price is blank for an OPEN trade, stored as an explicit empty
    # text value (not simply cleared).
    $ws.Range("G$row").Value = "'"
    $ws.Range("G$row").Style = "Normal"

    $ws.Range("H$row").Value = "OPEN"
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 100.7387999840491
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"

    # Exit Reason is blank too (trade still OPEN).
    $ws.Range("P$row").Value = "'"
    $ws.Range("P$row").Style = "Normal"

    $ws.Range("Q$row").Value = 0
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades 40

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking 40
